# #5: property boat&car done
# Fixes the "汽車" (car) sheet: row 1 was a stray duplicate of row 2's data
# instead of the proper column headers, and the per-row metadata columns
# (property_category .. index) that every other sheet already carries were
# missing. This restores the header row and appends those columns to the
# two existing data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("汽車")

# --- Row 1: replace the bogus duplicated-data header with real column names ---
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "capacity"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "register_date"
$ws.Range("F1").Value = "register_reason"
$ws.Range("G1").Value = "acquire_value"
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# --- Row 2 (record #32): append the source-tracking metadata columns ---
$ws.Range("H2").Value = "land"
$ws.Range("I2").Value = "normal"
$ws.Range("J2").Value = "2012-04-18"
$ws.Range("K2").Value = "蔡正元"
$ws.Range("L2").Value = 966
$ws.Range("M2").Value = "tmp671f1"
$ws.Range("N2").Value = 32

# --- Row 3 (record #33): same metadata columns ---
$ws.Range("H3").Value = "land"
$ws.Range("I3").Value = "normal"
$ws.Range("J3").Value = "2012-04-18"
$ws.Range("K3").Value = "蔡正元"
$ws.Range("L3").Value = 966
$ws.Range("M3").Value = "tmp671f1"
$ws.Range("N3").Value = 33
